$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.690.66"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "3.774.47"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.72"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.92"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "3.771.84"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.30"
$ws.Range("E11").Value = "  -1.75%  "
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("E13").Value = "  -2.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.96"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "4.407.54"
$ws.Range("D16").Value = "3.786.27"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").Value = "67.677.89"
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.38"
$ws.Range("E18").Value = "  +2.43%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.03"
$ws.Range("E21").Value = "  -6.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "456.96"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("E24").Value = "  +1.81%  "
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.95"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  -2.82%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.02"
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.24"
$ws.Range("E31").Value = "  +2.75%  "
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.64"
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").Value = "3.726.79"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0998"
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("E38").Value = "  -1.80%  "
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("E41").Value = "  -0.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.12"
$ws.Range("E44").Value = "  +2.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.23"
$ws.Range("E45").Value = "  +3.00%  "
$ws.Range("E46").Value = "  -1.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "149.51"
$ws.Range("E47").Value = "  +2.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.31"
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "390.01"
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("E50").Value = "  -5.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.23"
$ws.Range("E51").Value = "  +0.81%  "
